# Generate Report for Handback
# -----------------------------------------------------------------------------
# This script reproduces the "Generate Report for Handback" commit:
#  - Updates the localization status text from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it appears.
#  - Fills in the "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns on the zh-cn and de-de sheets for
#    both data rows, including brand new hyperlinks to a.md.
#  - Widens a handful of columns that now hold longer text.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$handedBackText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9b89aa8d1275e6929011ecfdd89e77d53d869f6/e2e/a.md"

# Hyperlink font look-alike (matches the workbook's existing "HyperLink" cell
# style: single underline, color FF6495ED).
$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

# -----------------------------------------------------------------------------
# 1. Overview sheet: status text + widened Priority/Content-Duplicate columns.
# -----------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $handedBackText
$wsOverview.Range("F2").Value = $handedBackText
$wsOverview.Range("E3").Value = $handedBackText
$wsOverview.Range("F3").Value = $handedBackText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# -----------------------------------------------------------------------------
# 2. zh-cn sheet
# -----------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $handedBackText
$wsZhCn.Range("C3").Value = $handedBackText

$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Range("I2").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor

$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-03 14:40:27"

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Range("I3").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor

$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-03 14:40:27"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# -----------------------------------------------------------------------------
# 3. de-de sheet
# -----------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $handedBackText
$wsDeDe.Range("C3").Value = $handedBackText

$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Range("I2").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 14:40:34"

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Range("I3").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor

$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 14:40:34"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
